# Revisão das regras de negócio
# Updates the "Regras de negócio" table on the "Prévia2" sheet:
#  - rewrites several business-rule descriptions (requisição de serviço -> leilão)
#  - drops the two stray comments that lived in column C (Colunas1)
#  - removes the now-empty RN36 row, shrinking the table from A1:C37 to A1:C36
#  - clears the yellow highlight fill that used to mark B3/B13/B16

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# --- Text updates in column B -------------------------------------------------
$ws.Range("B23").Value = "O Autônomo pode fazer até 3 ofertas por dia em qualquer serviço proposto."
$ws.Range("B26").Value = "O cliente pode abrir leilões ilimitadamente."
$ws.Range("B27").Value = "Para o cliente abrir um leilão, é necessário ser cadastrado no sistema."
$ws.Range("B28").Value = "O leilão poderá ficar aberto por até 30 dias."
$ws.Range("B31").Value = "O leilão deve conter data de abertura e de término, data de previsão, título e descrição do serviço. "
$ws.Range("B32").Value = "O cliente será notificado de novos lances."
$ws.Range("B33").Value = "O autônomo poderá fazer pesquisa de leilões."
$ws.Range("B34").Value = "Os lances mais novos do leilão devem ser exibidos primeiro."
$ws.Range("B35").Value = "Um autonomo pode procurar novos leilões e visualizá-los."
$ws.Range("B36").Value = "Autonomo não pode fazer lances em areas de atuação que não sejam pertencentes a ele."

# --- Drop the stray column-C comments ----------------------------------------
$ws.Range("C23").ClearContents()
$ws.Range("C32").ClearContents()

# --- The old RN36 row is gone entirely; delete it and let the table shrink ---
$ws.Rows.Item(37).Delete()

# --- Clear the yellow highlight fill on the cells that used to carry it ------
$ws.Range("B3").Interior.ThemeColor = 2
$ws.Range("B13").Interior.ThemeColor = 2
$ws.Range("B16").Interior.ThemeColor = 2

# --- Match the author's final view/selection state ---------------------------
$ws.Activate()
$ws.Range("B39").Select()
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
